$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D; this shifts existing data D:K -> E:L
$ws.Columns("D").Insert()

# Copy cell formatting (number format / style) from column E into the new column D
# (done per contiguous data block so we do not create stray cells in separator rows)
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# Populate the new column D with the FY2018 figures
$ws.Range("D7").Value = 43343
$ws.Range("D8").Value = 193300
$ws.Range("D9").Value = 138300
$ws.Range("D10").Value = 55000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 183000
$ws.Range("D18").Value = 10300
$ws.Range("D20").Value = -200
$ws.Range("D21").Value = 11100
$ws.Range("D22").Value = "NA"
$ws.Range("D23").Value = 10100
$ws.Range("D24").Value = 3200
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 6900
$ws.Range("D27").Value = 6900
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 200
$ws.Range("D33").Value = 6900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 6900
$ws.Range("D38").Value = 43343
$ws.Range("D41").Value = 2700
$ws.Range("D42").Value = 2800
$ws.Range("D43").Value = 26300
$ws.Range("D44").Value = 30500
$ws.Range("D45").Value = 2500
$ws.Range("D46").Value = 64900
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 9800
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 1600
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 76300
$ws.Range("D57").Value = 17700
$ws.Range("D58").Value = 100
$ws.Range("D59").Value = 8400
$ws.Range("D60").Value = 26200
$ws.Range("D61").Value = 8200
$ws.Range("D62").Value = "NA"
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 34400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 28500
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 41900
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43343
$ws.Range("D81").Value = 6900
$ws.Range("D83").Value = 1000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 5000
$ws.Range("D91").Value = -1600
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -2600
$ws.Range("D96").Value = -100
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -3700
$ws.Range("D101").Value = 200
$ws.Range("D102").Value = -1200
